# Cross browser Testing and Validations to test cases
# Add registration-form columns (firstname/lastname/regemail/regpwd) and a
# second "confirm email" row, mirroring the FB_TestData.xlsx commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (C1:F1) ------------------------------------------------
$ws.Range("C1").Value = "firstname"
$ws.Range("D1").Value = "lastname"
$ws.Range("E1").Value = "regemail"

# --- New data row (C2:F2) ---------------------------------------------------
$ws.Range("C2").Value = "sou"
$ws.Range("D2").Value = "muthu"
$ws.Range("E2").Value = "venkatshamuthu@gmail.com"

# F column (regpwd header + sample numeric pwd) -- added after E so the
# shared-string table fills up in the same order as the source workbook.
$ws.Range("F1").Value = "regpwd"
$ws.Range("F2").Value = 1234

# --- Extra confirmation row reusing the same e-mail as A2 -------------------
$ws.Range("A3").Value = "venkatsoumuthu@gmail.com"

# --- Hyperlinks for the two newly added e-mail cells ------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:venkatshamuthu@gmail.com")
$ws.Range("E2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:venkatsoumuthu@gmail.com")
$ws.Range("A3").Style = "Hyperlink"

# --- Column E width (matches the ~27.57-char wide column in the commit) ----
$ws.Columns.Item(5).ColumnWidth = 26.67

# --- Selection moves to B3 after the edits ----------------------------------
$ws.Range("B3").Select()
